$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.350.16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.731.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.38"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +13.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.267"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0636"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0898"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.974.85"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.730.49"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.29"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.561"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.78"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "28.308.91"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "243.62"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.66"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.72"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.55"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.74"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.85%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.21"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.14%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.507.47"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.70%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.21%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.970"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.608"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.40"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.92"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.879.36"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.806"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.73"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +9.94%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0114"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.64%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "91.06"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.21"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.08%  "
